# Adds the $coleccionVocales, $coleccionConsonant1 and $coleccionConsonant2
# example data structures to the worksheet, mirroring the existing
# $coleccionPalabras / $coleccionPartidas examples above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Style source cells (already present in the sheet):
#   B3   -> numeric "index" header style (bold, centered)      -> style 2
#   B4   -> value style (centered, thin border box)            -> style 1
#   B6   -> "Informacion de la estructura:" label style (bold) -> style 3
# ---------------------------------------------------------------------

# =========================================================
# $coleccionVocales  (rows 29-30, columns B..F -> 5 vowels)
# =========================================================
$ws.Range("A29").Value = "`$coleccionVocales"

$ws.Range("B3:F3").Copy()
$ws.Range("B29:F29").PasteSpecial(-4122)
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 3
$ws.Range("F29").Value = 4

$ws.Range("B4:F4").Copy()
$ws.Range("B30:F30").PasteSpecial(-4122)
$ws.Range("B30").Value = """A"""
$ws.Range("C30").Value = """E"""
$ws.Range("D30").Value = """I"""
$ws.Range("E30").Value = """O"""
$ws.Range("F30").Value = """U"""
$excel.CutCopyMode = $false

$ws.Range("B6").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B32").Value = "Información de la estructura:"

$ws.Range("B33").Value = "Tipo: Indexado (los índices son numéricos)"
$ws.Range("B34").Value = "Tipos de datos: Almacena valores String"
$ws.Range("B35").Value = "¿Para qué se utiliza?: guardar las vocales para ir calculando el puntaje en las partidas"

# ======================================================================
# $coleccionConsonant1  (rows 38-39, columns B..K -> 10 consonants B..M)
# ======================================================================
$ws.Range("A38").Value = "`$coleccionConsonant1"

$ws.Range("B3:F3").Copy()
$ws.Range("B38:K38").PasteSpecial(-4122)
$ws.Range("B38").Value = 0
$ws.Range("C38").Value = 1
$ws.Range("D38").Value = 2
$ws.Range("E38").Value = 3
$ws.Range("F38").Value = 4
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 6
$ws.Range("I38").Value = 7
$ws.Range("J38").Value = 8
$ws.Range("K38").Value = 9

$ws.Range("B4:F4").Copy()
$ws.Range("B39:K39").PasteSpecial(-4122)
$ws.Range("B39").Value = """B"""
$ws.Range("C39").Value = """C"""
$ws.Range("D39").Value = """D"""
$ws.Range("E39").Value = """F"""
$ws.Range("F39").Value = """G"""
$ws.Range("G39").Value = """H"""
$ws.Range("H39").Value = """J"""
$ws.Range("I39").Value = """K"""
$ws.Range("J39").Value = """L"""
$ws.Range("K39").Value = """M"""
$excel.CutCopyMode = $false

$ws.Range("B6").Copy()
$ws.Range("B41").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B41").Value = "Información de la estructura:"

$ws.Range("B42").Value = "Tipo: Indexado (los índices son numéricos)"
$ws.Range("B43").Value = "Tipos de datos: Almacena valores String"
$ws.Range("B44").Value = "¿Para qué se utiliza?: guardar las consonantes desde la B hasta la M para ir calculando el puntaje en las partidas"

# ======================================================================
# $coleccionConsonant2  (rows 47-48, columns B..L -> 11 consonants N..Z)
# ======================================================================
$ws.Range("A47").Value = "`$coleccionConsonant2"

$ws.Range("B3:F3").Copy()
$ws.Range("B47:L47").PasteSpecial(-4122)
$ws.Range("B47").Value = 0
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = 2
$ws.Range("E47").Value = 3
$ws.Range("F47").Value = 4
$ws.Range("G47").Value = 5
$ws.Range("H47").Value = 6
$ws.Range("I47").Value = 7
$ws.Range("J47").Value = 8
$ws.Range("K47").Value = 9
$ws.Range("L47").Value = 10

$ws.Range("B4:F4").Copy()
$ws.Range("B48:L48").PasteSpecial(-4122)
$ws.Range("B48").Value = """N"""
$ws.Range("C48").Value = """P"""
$ws.Range("D48").Value = """Q"""
$ws.Range("E48").Value = """R"""
$ws.Range("F48").Value = """S"""
$ws.Range("G48").Value = """T"""
$ws.Range("H48").Value = """V"""
$ws.Range("I48").Value = """W"""
$ws.Range("J48").Value = """X"""
$ws.Range("K48").Value = """Y"""
$ws.Range("L48").Value = """Z"""
$excel.CutCopyMode = $false

$ws.Range("B6").Copy()
$ws.Range("B50").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B50").Value = "Información de la estructura:"

$ws.Range("B51").Value = "Tipo: Indexado (los índices son numéricos)"
$ws.Range("B52").Value = "Tipos de datos: Almacena valores String"
$ws.Range("B53").Value = "¿Para qué se utiliza?: guardar las consonantes desde la N hasta la Z para ir calculando el puntaje en las partidas"

# ---------------------------------------------------------------------
# Column A width & view state, matching the author's final edit.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.2857142857143

$ws.Range("D57").Select()
